$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for B2:B452 (row index 0..450 in column A corresponds to these)
$values = @(
    2139147823.158982, 2140169670.86859, 2143006225.766772, 2145082632.953526, 2147250130.976365, 2149502752.228405, 2151836197.010477, 2154246076.736089,
    2156727943.340976, 2159279673.328859, 2161897905.419241, 2164580232.194396, 2167325258.146171, 2170130885.475006, 2172996456.420834, 2175919923.288563,
    2178899825.238841, 2181936479.223039, 2185029097.037541, 2188177615.959925, 2191381254.212604, 2194640043.447813, 2197953359.703855, 2201321691.547487,
    2204743578.236913, 2208221349.623917, 2211753354.391282, 2215341290.482126, 2218985171.457137, 2222684564.047967, 2226440416.760038, 2230253008.195704,
    2234122640.283274, 2238050664.409834, 2242037246.023265, 2246084143.734515, 2250189602.445416, 2254356201.511071, 2258582987.702851, 2262872827.979146,
    2267225023.366434, 2271641614.253325, 2276121644.582702, 2280667047.313608, 2285277783.370313, 2289955927.032278, 2294702488.236923, 2299516962.239019,
    2304401150.112209, 2309355054.465537, 2314381081.591636, 2319480237.381257, 2324652189.052814, 2329897482.322262, 2335218689.794785, 2340617252.448229,
    2346092446.809124, 2351645740.200575, 2357280106.534969, 2362994391.860372, 2368790237.764625, 2374669453.909499, 2380633809.091628, 2386682183.258844,
    2392817442.717723, 2399041524.821189, 2405352748.545685, 2411754975.722863, 2418249635.312344, 2424835226.815031, 2431516855.364607, 2438291933.350747,
    2445165028.537382, 2452136243.614618, 2459205981.293279, 2466376941.494812, 2473649471.820322, 2481026812.74739, 2488508204.51603, 2496097064.09779,
    2503793301.476742, 2511599375.86972, 2519515456.602367, 2527545165.693846, 2535688243.738586, 2543948415.724522, 2552324213.241264, 2560820366.346466,
    2569436219.580472, 2578175422.523264, 2587037059.409885, 2596025283.255594, 2605140860.417963, 2614384911.100444, 2623761742.671799, 2633269355.035767,
    2642911741.938412, 2652691338.678057, 2662609129.544691, 2672666872.084159, 2682867561.085286, 2693210949.473764, 2703700965.127103, 2714340135.985356,
    2725127419.452473, 2736068690.234909, 2747163776.903384, 2758415356.807035, 2769824332.483077, 2781394926.24462, 2793129287.573997, 2805027453.141021,
    2817093177.915387, 2829330363.114701, 2841738505.216976, 2854319420.944006, 2867078130.00293, 2880017028.723805, 2893137214.586604, 2906440284.929407,
    2919930314.606479, 2933609317.310183, 2947481751.874609, 2961545193.459412, 2975808765.784482, 2990270576.353929, 3004936037.627391, 3019804339.507375,
    3034882721.764188, 3050170808.036109, 3065673582.096919, 3081390739.207226, 3097330023.65723, 3113490313.269359, 3129876895.341556, 3146490444.521179,
    3163338632.673133, 3180418515.318353, 3197738296.479166, 3215298817.577816, 3233104900.385302, 3251155002.896908, 3269460619.719211, 3288018668.590996,
    3306835148.888885, 3325912566.108235, 3345257882.163918, 3364868163.201485, 3384754142.00019, 3404913971.948538, 3425354347.122975, 3446078012.551781,
    3467089727.462876, 3488392103.980564, 3509991785.26037, 3531887713.417586, 3554087344.578665, 3576593357.562255, 3599413495.166563, 3622548704.787588,
    3646004608.695323, 3669784443.558047, 3693894282.068661, 3718338316.53479, 3743119363.711652, 3768242997.329861, 3793716094.763052, 3819537858.50223,
    3845720261.252741, 3872262011.057324, 3899173406.190725, 3926452633.069948, 3954113633.08545, 3982152014.401069, 4010582178.583288, 4039400533.49576,
    4068622488.355793, 4098242183.004927, 4128276109.423006, 4158720148.646173, 4189589396.979265, 4220879997.661005, 4252608596.335671, 4284769598.618806,
    4317380280.465104, 4350438009.597613, 4383955370.067511, 4417932471.294792, 4452383446.90184, 4487306752.007385, 4522716969.941672, 4558613113.214097,
    4595009067.086094, 4631905593.780021, 4669316142.81266, 4707243232.204408, 4745695000.872828, 4784680294.932563, 4824207020.432088, 4864279680.89009,
    4904910425.515305, 4946105363.485628, 4987871031.267405, 5030217847.187349, 5073153490.236794, 5116685994.207702, 5160825302.426821, 5205579525.413706,
    5250957804.676571, 5296967493.115937, 5343620616.331939, 5390925842.201008, 5438892058.614079, 5487530353.735115, 5536849654.539299, 5586859987.328229,
    5637573701.623627, 5688999604.74789, 5741148300.531031, 5794031852.236032, 5847662355.055492, 5902047560.65594, 5957202150.803892, 6013137729.137671,
    6069864712.114258, 6127397994.626786, 6185748254.494116, 6244930060.184892, 6304950659.69581, 6365833369.551651, 6427584360.373694, 6490219036.532533,
    6553751760.188979, 6618200210.280345, 6683572451.112101, 6749888784.088862, 6817164677.950579, 6885411163.120297, 6954650393.382157, 7024895264.294799,
    7096161657.424824, 7168466632.269297, 7241834033.193437, 7316271516.063814, 7391806015.266739, 7468450212.945902, 7546228122.995893, 7625154213.176849,
    7705251968.983586, 7786540301.734158, 7869037514.177157, 7952771617.135534, 8037756426.471581, 8124019399.261948, 8211580561.644727, 8300466206.80072,
    8390693222.602459, 8482294327.741354, 8575288046.972179, 8669702688.520958, 8765561844.218019, 8862892535.830854, 8961724400.23304, 9062079593.599874,
    9163993849.080465, 9267487728.877844, 9372596864.651907, 9479348595.300964, 9587776733.461727, 9697906982.918606, 9809778944.442406, 9923418144.536459,
    10038862758.16077, 10156145475.21699, 10275300547.60483, 10396364474.53922, 10519376783.96411, 10644371013.32281, 10771384345.86182, 10900462157.8933,
    11031634818.79165, 11164951700.27343, 11300447395.31211, 11438166220.82712, 11578154504.23804, 11720451551.56901, 11865103897.39263, 12012155083.45255,
    12161656207.08433, 12313648948.72556, 12468189206.61028, 12625316172.59946, 12785088233.46329, 12947550588.31529, 13112758098.04309, 13280763176.77412,
    13451616755.06759, 13625374929.47569, 13802094465.56626, 13981826297.79817, 14164633844.13758, 14316868065.51152, 14471762976.59959, 14629356423.62662,
    14789701432.86103, 14952854860.02998, 15118864149.65591, 15287785962.24028, 15459667543.51037, 15634567661.29946, 15812537768.81039, 15993633965.94591,
    16177912406.19309, 16365424963.65233, 16556227241.97239, 16750383280.32377, 16947939936.41692, 17148960817.36503, 17353499841.90049, 17561611466.70222,
    17773357356.60613, 17988794039.65254, 18207977871.21556, 18430967124.36834, 18657815125.27998, 18888582460.21159, 19123324802.98704, 19362098302.34279,
    19604958526.59029, 19851960769.97374, 20103160256.06564, 20358608032.68411, 20618359855.12103, 20882468914.92345, 21150983720.67014, 21423959900.68146,
    21701440064.46295, 21983474670.58469, 22270111632.15013, 22561396653.60631, 22857369628.1791, 23158076113.15362, 23463558484.33343, 23773848557.26277,
    24088987485.42239, 24409008161.41386, 24733944695.76267, 25063823240.25132, 25398677551.06616, 25738525058.45218, 26083393960.61514, 26433298864.6455,
    26788263756.62401, 27148296309.24984, 27513408560.35504, 27883609269.08095, 28252237678.00124, 28625875779.3116, 29004511842.48681, 29388149000.37703,
    29776774629.93819, 30170376247.73244, 30568938538.17593, 30972436988.19038, 31380853682.52018, 31794156321.82282, 32212315888.89773, 32635294810.46419,
    33063055342.58815, 33495551448.69551, 33932740225.86945, 34311276449.39563, 34693710891.51314, 35079985545.53459, 35470042326.07714, 35863820140.10318,
    36261253916.80518, 36661835025.06287, 37060657150.15664, 37445280252.47012, 37801221174.68129, 38281363055.73478, 38722099577.17403, 39127595996.00301,
    39505964630.07344, 39866837070.31576, 40219670336.78915, 40572750673.22153, 40932781465.53497, 41304816810.4586, 41692402273.93326, 41827717094.72745,
    41981831130.44463, 42154945303.14653, 42346676125.11169, 42556235806.90793, 42782591116.08569, 43024549460.70847, 43280858028.25491, 43550246105.06527,
    43831476056.01742, 43846682916.69071, 43870435216.3786, 43901577399.4832, 43939050708.56857, 43981896436.82764, 44029244265.40927, 44079664009.32339,
    44125313046.30697, 44147695473.47038, 44126068945.64772, 44129052836.92637, 44072392531.87798, 43963066398.10667, 43813332186.61668, 43637204936.91028,
    43448042404.3896, 43257244575.64023, 43073755395.8889, 42904079024.22314, 42752542122.01902, 42621688893.77557, 42512655695.48553, 42425542108.17752,
    42359700284.71114, 42313984424.3711, 42286930749.05841, 42276900925.53999, 42282196072.48995, 42301102752.0619, 42331968495.3241, 42373217783.36245,
    42423374507.5332, 42481065636.46476, 42545033129.1634, 42614125445.29008, 42687303672.4599, 42763621336.10094, 42842223439.48396, 42922340495.83239,
    43003280422.73766, 43084421225.01428, 43165200877.45817, 43245123608.50325, 43323738424.61771, 43400639671.84509, 43475464437.91836, 43547893928.89145,
    43617636006.85989, 43684431855.56743, 380041927.7425824
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($i + 2, 2).Value = $values[$i]
}
